$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" everywhere it appears ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- Narrow the "zh-cn"/"de-de" status columns ---
# Overview sheet: columns E (zh-cn) and F (de-de)
$overview.Columns.Item(5).ColumnWidth = 12.58
$overview.Columns.Item(6).ColumnWidth = 12.58

# zh-cn sheet: column C (Status)
$zhcn.Columns.Item(3).ColumnWidth = 12.58

# de-de sheet: column C (Status)
$dede.Columns.Item(3).ColumnWidth = 12.58
